$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 181; this shifts the existing rows 181-225
# down to 182-226 (and extends the used range to A1:R226), matching the
# diff exactly for the rows that simply move down one position.
$ws.Rows.Item(181).Insert()

# Populate the freshly inserted row 181 with the new weekly record.
$ws.Range("A181").Value = 4
$ws.Range("B181").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C181").Value = "Los Lagos"
$ws.Range("D181").Value = 44641
$ws.Range("E181").Value = 10
$ws.Range("F181").Value = 100112044
$ws.Range("G181").Value = "Perejil"
$ws.Range("H181").Value = "Sin especificar"
$ws.Range("I181").Value = "Primera"
$ws.Range("J181").Value = 70
$ws.Range("K181").Value = 5000
$ws.Range("L181").Value = 5000
$ws.Range("M181").Value = 5000
$ws.Range("N181").Value = '$/docena de atados (2 kilos)'
$ws.Range("O181").Value = 'Región de La Araucanía'
$ws.Range("P181").Value = 2500
$ws.Range("Q181").Value = 2
$ws.Range("R181").Value = "Hortaliza"
